$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Date - 26.12.2023"
$ws.Range("D4").Value = "DEL.N1LN3L.WD.D.GPD.DMP : last exported on 23.12.23`nDEL.N1LN3L.WD.D.GPL.DMP : last exported on 23.12.23`nDEL.N1LN3L.WD.D.GTR.DMP : last exported on 26.12.23`nDEL.N1LKQQ.WD.D.UPD.DMP : last exported on 23.12.23`nDEL.N1LN3L.WD.D.ASPD.DMP : last exported on 23.12.23`nDEL.N1LN3L.WD.D.ASPL.DMP : last exported on 23.12.23`nDEL.N1LN0R.WD.D.CPD.DMP : last exported on 23.12.23`nDEL.N1LN0R.WD.D.CPL.DMP : last exported on 23.12.23`nDEL.KT6N5M.SL.D.GGO.ZIP : last exported on 23.12.23`nDEL.KT6N5M.SL.D.UGO.ZIP : last exported on 23.12.23`nDEL.KT6N5M.SL.D.ASGO.ZIP : last exported on 23.12.23"
$ws.Range("D7").Value = "DEL.N1LN3L.WI.D.GPD.DMP : last exported on 26.12.23`nDEL.N1LN3L.WI.D.GPL.DMP : last exported on 26.12.23`nDEL.N1LN3L.WI.D.GTR.DMP : last exported on 26.12.23`nDEL.N1LN3L.WI.D.ASPD.DMP : last exported on 26.12.23`nDEL.N1LN3L.WI.D.ASPL.DMP : last exported on 26.12.23`nDEL.N1LN0R.WI.D.CPD.DMP : last exported on 26.12.23`nDEL.N1LN0R.WI.D.CPL.DMP : last exported on 26.12.23`nDEL.KT6E35.WI.GGO.ZIP : last exported on 25.12.23`nDEL.KT6E35.WI.IMG.GGO.ZIP : last exported on 25.12.23`nDEL.KT6N5M.WI.GGO.ZIP : last exported on 25.12.23`nDEL.KT6N5M.WI.IMG.GGO.ZIP : last exported on 25.12.23"
$ws.Range("D10").Value = "DEL.N1LN3L.TP.D.GPD.DMP : last exported on 25.12.23`nDEL.N1LN3L.TP.D.GPL.DMP : last exported on 25.12.23`nDEL.N1LN3L.TP.D.GTR.DMP : last exported on 25.12.23`nDEL.N1LKQQ.TP.D.UPD.DMP : last exported on 25.12.23`nDEL.N1LN3L.TP.D.ASPD.DMP : last exported on 25.12.23`nDEL.N1LN3L.TP.D.ASPL.DMP : last exported on 25.12.23`nDEL.N1LN0R.TP.D.CPD.DMP : last exported on 25.12.23`nDEL.N1LN0R.TP.D.CPL.DMP : last exported on 25.12.23`nDEL.KT6E35.TP.D.GGO.ZIP : last exported on 24.12.23`nDEL.KT6E35.TP.D.UGO.ZIP : last exported on 24.12.23`nDEL.KT6N5M.TP.D.GGO.ZIP : last exported on 24.12.23`nDEL.KT6N5M.TP.D.UGO.ZIP : last exported on 24.12.23`nDEL.KT6N5M.TP.D.ASGO.ZIP : last exported on 24.12.23`nDEL.KT6N5M.TP.D.UOC.ZIP : last exported on 24.12.23`nDEL.KT6N5M.TP.D.ASOC.ZIP : last exported on 24.12.23`nDEL.KT6N5M.TP.D.GOC.ZIP : last exported on 24.12.23`nDED.KT6I1M.HST2.D.GGO.ZIP : last exported on 24.12.23"
$ws.Range("D11").Value = "DEL.N1LN3L.MC.F.GPD.DMP : last exported on 30.11.23`nDEL.N1LN3L.MC.F.GPL.DMP : last exported on 30.11.23`nDEL.N1LN3L.MC.F.GTR.DMP : last exported on 25.12.23`nDEL.N1LKQQ.MC.F.UPD.DMP : last exported on 30.11.23`nDEL.N1LN3L.MC.F.ASPD.DMP : last exported on 30.11.23`nDEL.N1LN3L.MC.F.ASPL.DMP : last exported on 30.11.23`nDGI.KT6R11.MANDANT.TXT : last exported on 30.11.23"
$ws.Range("D13").Value = "DEL.N1LN3L.AU.D.GPD.DMP : last exported on 23.12.23`nDEL.N1LN3L.AU.D.GPL.DMP : last exported on 23.12.23`nDEL.N1LN3L.AU.D.GTR.DMP : last exported on 25.12.23`nDEL.N1LN0R.AU.D.CPD.DMP : last exported on 23.12.23`nDEL.N1LN0R.AU.D.CPL.DMP : last exported on 23.12.23`nDEL.KT6E35.AU.GGO.ZIP : last exported on 23.12.23`nDEL.KT6E35.AU.IMG.GGO.ZIP : last exported on 23.12.23"
$ws.Range("D15").Value = "DGI.N1LR11.FB.ORGUNIT.CSV : last exported on 25.12.23`nDGI.N1LR11.FB.IMPTR.CSV : last exported on 25.12.23"
$ws.Range("D16").Value = "DEL.N1LN3L.DP.F.GPD.DMP : last exported on 25.12.23`nDEL.N1LN3L.DP.F.GPL.DMP : last exported on 25.12.23`nDEL.N1LN3L.DP.F.GTR.DMP : last exported on 25.12.23`nDEL.N1LKQQ.DP.F.UPD.DMP : last exported on 25.12.23`nDEL.N1LN3L.DP.F.ASPD.DMP : last exported on 25.12.23`nDEL.N1LN3L.DP.F.ASPL.DMP : last exported on 25.12.23`nDEL.N1LN0R.DP.F.CPD.DMP : last exported on 25.12.23`nDEL.N1LN0R.DP.F.CPL.DMP : last exported on 25.12.23`nDEL.N1LN3L.DISSORG.ZIP : last exported on 25.12.23`nDEL.N1LN3L.DISSAORG.ZIP : last exported on 25.12.23"
$ws.Range("D19").Value = "DEL.N1LN3L.MT.F.GPD.DMP : last exported on 25.12.23`nDEL.N1LN3L.MT.F.GPL.DMP : last exported on 25.12.23`nDEL.N1LN3L.MT.F.GTR.DMP : last exported on 25.12.23`nDEL.N1LKQQ.MT.F.UPD.DMP : last exported on 25.12.23`nDEL.N1LN3L.MT.F.ASPD.DMP : last exported on 25.12.23`nDEL.N1LN3L.MT.F.ASPL.DMP : last exported on 25.12.23`nDEL.N1LN0R.MT.F.CPD.DMP : last exported on 25.12.23`nDEL.N1LN0R.MT.F.CPL.DMP : last exported on 25.12.23`nDEL.N1LR31.MT.F.GPD.DMP : last exported on 25.12.23`nDEL.N1LI4X.MT.F.GPD.ZIP : last exported on 25.12.23`nDEL.N1LI4Y.MT.F.GPL.ZIP : last exported on 25.12.23`nDEL.N1LGN1.NL.F.GPD.ZIP : last exported on 25.12.23`nDEL.N1LKT6.MT.??????.G : last exported on 25.12.23`nDEL.N1LKT6.MT.??????.U : last exported on 25.12.23`nDEL.N1LKT6.MT.??????.AS : last exported on 25.12.23`nDVL.N1LN5X.VLM.DSP : last exported on 25.12.23`nDVL.N1LN5X.VLM.WT : last exported on 25.12.23`nDDS.N1LR11.DSP : last exported on 25.12.23`nDDS.N1LR11.WT : last exported on 25.12.23`nDEL.KT6N5M.MT.F.GGO.ZIP : last exported on 23.12.23`nDEL.KT6N5M.MT.F.UGO.ZIP : last exported on 23.12.23`nDEL.KT6N5M.MT.F.ASGO.ZIP : last exported on 23.12.23`nDEL.KT6E35.MT.F.GGO.ZIP : last exported on 23.12.23`nDEL.KT6E35.MT.F.UGO.ZIP : last exported on 23.12.23`nDEL.KT6N5X.MT.F.GGO.ZIP : last exported on 23.12.23`nDEL.KT6N5X.MT.F.UGO.ZIP : last exported on 23.12.23`nDEL.KT6N5X.MT.F.ASGO.ZIP : last exported on 23.12.23"
$ws.Range("D20").Value = "DEL.N1LN3L.EL.F.GPD.DMP : last exported on 25.12.23`nDEL.N1LN3L.EL.F.GPL.DMP : last exported on 25.12.23`nDEL.N1LE35.ELFI.G.DATA : last exported on 25.12.23`nDEL.N1LN5M.EL.F.GGO.ZIP : last exported on 25.12.23`nDEL.N1LN5X.EL.F.GGO.ZIP : last exported on 25.12.23`nDEL.N1LE35.EL.F.GGO.ZIP : last exported on 25.12.23`nDEL.N1LKT6.EC.??????.G : last exported on 25.12.23`nDWT.N1LN3L.WT.??????.ZIP : last exported on 25.12.23`nDEL.N1LKQQ.EL.F.UPD.DMP : last exported on 25.12.23`nDEL.N1LE35.ELFI.U.DATA : last exported on 25.12.23`nDEL.N1LN5M.EL.F.UGO.ZIP : last exported on 25.12.23`nDEL.N1LN5X.EL.F.UGO.ZIP : last exported on 25.12.23`nDEL.N1LE35.EL.F.UGO.ZIP : last exported on 25.12.23`nDEL.N1LKT6.EC.??????.U : last exported on 25.12.23`nDEL.N1LN3L.EL.F.ASPD.DMP : last exported on 25.12.23`nDEL.N1LN3L.EL.F.ASPL.DMP : last exported on 25.12.23`nDEL.N1LN0R.EL.F.CPD.DMP : last exported on 25.12.23`nDEL.N1LN0R.EL.F.CPL.DMP : last exported on 25.12.23`nDEL.N1LN5M.EL.F.ASGO.ZIP : last exported on 25.12.23`nDEL.N1LN5X.EL.F.ASGO.ZIP : last exported on 25.12.23`nDEL.N1LE35.EL.F.ASGO.ZIP : last exported on 25.12.23`nDEL.N1LKT6.EC.??????.AS : last exported on 25.12.23`nDEL.N1LN3L.EL.F.GTR.DMP : last exported on 25.12.23`nDGV.N1LR31.EL.F.HST.ZIP : last exported on 25.12.23`nDGV.N1LN3L.EL.F.HST.ZIP : last exported on 25.12.23`nDGV.N1LCEN.EL.F.ZIP : last exported on 25.12.23`nDGV.N1LN7K.EL.F.ZIP : last exported on 25.12.23`nDGV.N1LN5X.EL.F.ZIP : last exported on 25.12.23`nDEL.KT6KT6.EL.F.DRE.CSV : last exported on 25.12.23`nDEL.KT6N5X.EL.F.AU.ZIP : last exported on 24.12.23"
